$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank row at position 9, pushing the old row 9 (total) and
#    row 10 (footer) down to rows 10 and 11.
$ws.Rows("9:9").Insert()

# 2) Copy the formatting + values of row 8 (the item-2 row, still holding the
#    "ERASTAPEX..." data at this point) into the newly inserted row 9. This
#    preserves cell styles/number formats exactly (same style ids as row 8).
$ws.Range("A8:Q8").Copy($ws.Range("A9:Q9"))

# 3) Fix row heights: new row 9 and the shifted total row (now row 10) both
#    end up at 25.5pt (row 10 was 26.25pt before the insert).
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 25.5

# 4) Row 9 is "item 3" now -- bump its item number from 2 to 3.
$ws.Range("A9").Value = 3

# 5) Re-create the merges for the new row 9 (Insert() does not duplicate the
#    merges from row 8 automatically).
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# 6) Now overwrite row 8 with the new item ("BRUFEN 400MG 30 TAB") data --
#    this is the item that got inserted ahead of "ERASTAPEX..." (now row 9).
$ws.Range("C8").Value = "BRUFEN 400MG 30 TAB"
$ws.Range("H8").Value = "1:1"
$ws.Range("N8").Value = "78.00"
$ws.Range("P8").Value = "25.7400"
$ws.Range("Q8").Value = "0:1"

# 7) Update the grand total cell (now row 10, was row 9): 165 + 25.74.
$ws.Range("N10").Value = 190.74000000000001
